# Add an "Other:" row to the SUPPLIES & RESOURCES section of the budget,
# just above the "SUPPLIES & RESOURCES TOTAL:" row (i.e. insert a new
# row 14, pushing everything from the old row 14 down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 14 (old row 14 "SUPPLIES & RESOURCES
# TOTAL:" and everything below shifts down to row 15+).
$ws.Rows("14:14").Insert()

# The row above (old "Lunch and snacks:" row, still row 13) has the same
# look as every other expense line: a plain label cell in A, a currency
# input in B, and a merged, wrap-text "Description" field across C:E.
# Re-create that merge on the new row first ...
$ws.Range("C14:E14").Merge()

# ... then copy that row's formatting (fonts/fills/borders/number
# formats) onto the freshly inserted row so the new line matches the
# rest of the table exactly.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

# Label the new row.
$ws.Range("A14").Value = "Other:"

# The "SUPPLIES & RESOURCES TOTAL:" row (now row 15) should sum through
# the newly added "Other:" row as well.
$ws.Range("B15").Formula = "=SUM(B9:B14)"

# Match the author's final selection.
$ws.Range("B14").Select()
